$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Paths"
$ws.Range("B2").Value = "Ejobs.xlsx"
$ws.Range("B3").Value = "Linkedin.xlsx"
$ws.Range("B4").Value = "BestJobs.xlsx"
$ws.Range("B5").Value = "AllJobs.xlsx"
$ws.Range("B6").Value = "NewJobs.xlsx"

$ws.Columns.Item(1).ColumnWidth = 72
$ws.Columns.Item(2).ColumnWidth = 13

$ws.Range("B7").Select()
